$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values remain plain text even when they look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.625.84'
$ws.Range("E2").Value = '  +2.35%  '
$ws.Range("D3").Value = '1.891.57'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '244.48'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4957'
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = '0.2958'
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("D9").Value = '0.06812'
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("D10").Value = '1.892.07'
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = '17.09'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").Value = '0.07317'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").Value = '91.08'
$ws.Range("E13").Value = '  +5.83%  '
$ws.Range("D14").Value = '5.094'
$ws.Range("E14").Value = '  +5.49%  '
$ws.Range("D15").Value = '0.6736'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").Value = '30.609.76'
$ws.Range("E16").Value = '  +2.30%  '
$ws.Range("D17").Value = '0.000007910'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '13.26'
$ws.Range("E19").Value = '  +4.60%  '
$ws.Range("D20").Value = '2.136.09'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '4.865'
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").Value = '179.72'
$ws.Range("E23").Value = '  +34.08%  '
$ws.Range("D24").Value = '6.048'
$ws.Range("E24").Value = '  +8.57%  '
$ws.Range("D25").Value = '9.297'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").Value = '154.30'
$ws.Range("E26").Value = '  +2.74%  '
$ws.Range("D27").Value = '18.80'
$ws.Range("E27").Value = '  +12.77%  '
$ws.Range("D28").Value = '1.926'
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("D29").Value = '1.387'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '4.323'
$ws.Range("E30").Value = '  +4.43%  '
$ws.Range("D31").Value = '0.08936'
$ws.Range("E31").Value = '  +3.09%  '
$ws.Range("D32").Value = '4.037'
$ws.Range("E32").Value = '  +2.87%  '
$ws.Range("D33").Value = '0.05200'
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").Value = '0.7369'
$ws.Range("E34").Value = '  +5.78%  '
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  +3.84%  '
$ws.Range("D36").Value = '2.680'
$ws.Range("D37").Value = '0.01868'
$ws.Range("E37").Value = '  +10.19%  '
$ws.Range("D38").Value = '2.697'
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("D39").Value = '2.167'
$ws.Range("E39").Value = '  +0.58%  '
$ws.Range("D40").Value = '0.9344'
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("D41").Value = '0.4356'
$ws.Range("E41").Value = '  +4.44%  '
$ws.Range("D42").Value = '105.95'
$ws.Range("E42").Value = '  +4.42%  '
$ws.Range("D43").Value = '5.810'
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +3.88%  '
$ws.Range("D46").Value = '0.1350'
$ws.Range("E46").Value = '  +7.99%  '
$ws.Range("D47").Value = '0.05841'
$ws.Range("E47").Value = '  +3.55%  '
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("D49").Value = '0.3891'
$ws.Range("E49").Value = '  +5.59%  '
$ws.Range("D50").Value = '8.491'
$ws.Range("E50").Value = '  +5.38%  '
$ws.Range("E51").Value = '  +3.90%  '
